$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "2025-05-01T11:14:47.873Z"
$ws.Range("B6").Value = "IDRF"
$ws.Range("C6").Value = "C3"
$ws.Range("D6").Value = "الرحلة 1"
$ws.Range("E6").Value = "الصمود"
$ws.Range("F6").Value = "يامن "
$ws.Range("G6").Value = "12"
$ws.Range("H6").Value = ""
